# Updates the cryptocurrency price/volume table on Sheet1 (row 1 = headers,
# rows 2-51 = coins) to match the latest coinranking.com snapshot.
#
# Most rows only get refreshed Price (D) / Volume(1h) (E) figures. A few rows
# (42-51) also changed which coin occupies that rank, so Coin (B) and Link (C)
# are rewritten there too -- including a brand-new entry, BabyDogeCoin, which
# bumped USDD off the bottom of the list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.572.66"
$ws.Range("E2").Value = "  +1.24%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.629.31"
$ws.Range("E3").Value = "  +1.07%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "

# Row 6: USDC
$ws.Range("E6").Value = "  +0.01%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.493"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.61%  "

# Row 8: Cardano
$ws.Range("E8").Value = "  +1.06%  "

# Row 9: Dogecoin
$ws.Range("E9").Value = "  +1.02%  "

# Row 10: Solana
$ws.Range("E10").Value = "  +1.14%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0838"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.31%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.855.94"
$ws.Range("E12").Value = "  +1.06%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.629.44"
$ws.Range("E13").Value = "  +1.07%  "

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.19%  "

# Row 15: Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.525"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.38%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "26.546.33"
$ws.Range("E16").Value = "  +1.13%  "

# Row 17: Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.81%  "

# Row 18: ShibaInu
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.27%  "

# Row 19: Dai
$ws.Range("E19").Value = "  +0.03%  "

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.66%  "

# Row 21: Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "

# Row 22: Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "

# Row 23: Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.11%  "

# Row 24: Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.72%  "

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.93%  "

# Row 26: BinanceUSD
$ws.Range("E26").Value = "  +0.03%  "

# Row 27: Stellar
$ws.Range("E27").Value = "  -0.29%  "

# Row 28: Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.36%  "

# Row 29: EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "

# Row 30: Hedera
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0515"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.30%  "

# Row 31: PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.96%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.08%  "

# Row 34: LidoDAOToken
$ws.Range("E34").Value = "  +0.22%  "

# Row 35: HuobiToken
$ws.Range("E35").Value = "  -0.33%  "

# Row 36: VeChain
$ws.Range("E36").Value = "  +4.06%  "

# Row 37: Maker
$ws.Range("D37").Value = "1.166.87"
$ws.Range("E37").Value = "  +0.76%  "

# Row 38: ARBITRUM
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.805"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.04%  "

# Row 39: PaxDollar
$ws.Range("E39").Value = "  +0.07%  "

# Row 40: ImmutableX
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.502"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.37%  "

# Row 41: MXToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.55%  "

# Row 42: TrustWalletToken -> FraxShare
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.35%  "

# Row 43: FraxShare -> TrustWalletToken
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.786"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "

# Row 44: RocketPoolETH
$ws.Range("D44").Value = "1.763.89"
$ws.Range("E44").Value = "  +0.95%  "

# Row 45: Quant
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "

# Row 46: RenderToken
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.34%  "

# Row 47: Aave -> BabyDogeCoin
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0103"
$ws.Range("E47").Value = "  -2.91%  "

# Row 48: Cronos -> Aave
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.15%  "

# Row 49: EnergySwap -> Cronos
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "

# Row 50: Mantle -> EnergySwap
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.21%  "

# Row 51: USDD -> Mantle
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.409"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.53%  "

Write-Output "Updated cryptos list"
